$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.848412752151489
$ws.Range("B1").Value = 2.03176736831665
$ws.Range("C1").Value = 2.397714376449585
$ws.Range("D1").Value = 2.996721982955933
$ws.Range("E1").Value = 3.043732881546021
